$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 'Вн'
$ws.Cells.Item(2, 6).Value = 'ДГИ-В-87772/24  14.08.2024'
$ws.Cells.Item(2, 7).Value = 'ДГИ-В-87772/24 14.08.2024'
$ws.Cells.Item(2, 8).Value = 'Кому: Мусиенко О.А. (Департамент городского имущества города Москвы)  От кого:  Мишиева Э.Ш. (Департамент городского имущества города Москвы)'
$ws.Cells.Item(2, 9).Value = 'УП. СВАО. Информация в отношении дома по адресу: ул. Молокова, д. 17-19'

# Row 3
$ws.Cells.Item(3, 5).Value = 'Вх'
$ws.Cells.Item(3, 6).Value = 'ДГИ-Э-129183/24  14.08.2024'
$ws.Cells.Item(3, 7).Value = '13/2 13.08.2024'
$ws.Cells.Item(3, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  --- ("-")'
$ws.Cells.Item(3, 9).Value = 'Запрос в отношении гр. Республики Беларусь Козодавенко А.В.'

# Row 4
$ws.Cells.Item(4, 5).Value = 'Вн'
$ws.Cells.Item(4, 6).Value = 'ДГИ-В-87616/24  14.08.2024'
$ws.Cells.Item(4, 7).Value = 'ДГИ-В-87616/24 14.08.2024'
$ws.Cells.Item(4, 8).Value = 'Кому: Мусиенко О.А. (Департамент городского имущества города Москвы)  От кого:  Демонова Л.В. (Департамент городского имущества города Москвы)'
$ws.Cells.Item(4, 9).Value = 'О показе жилых помещений ЮВАО'

# Row 5
$ws.Cells.Item(5, 5).Value = 'Гр'
$ws.Cells.Item(5, 6).Value = 'ДГИ-ЭГР-46949/24  14.08.2024'
$ws.Cells.Item(5, 7).Value = '56023306 14.08.2024'
$ws.Cells.Item(5, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Обращение граждан (Обращение граждан)'
$ws.Cells.Item(5, 9).Value = 'Обращения граждан Вопрос 1. Сообщение с mos.ru, идентификатор: 56023306 Исаев Александр Петрович Жалоба , жилищный учет'

# Row 6
$ws.Cells.Item(6, 5).Value = 'Вн'
$ws.Cells.Item(6, 6).Value = 'ДГИ-В-87613/24  14.08.2024'
$ws.Cells.Item(6, 7).Value = 'ДГИ-В-87613/24 14.08.2024'
$ws.Cells.Item(6, 8).Value = 'Кому: Мусиенко О.А. (Департамент городского имущества города Москвы)  От кого:  Демонова Л.В. (Департамент городского имущества города Москвы)'
$ws.Cells.Item(6, 9).Value = 'Об осмотре жилых помещений по КПИ и ДСН в ЮВАО'

# Row 7
$ws.Cells.Item(7, 5).Value = 'Вн'
$ws.Cells.Item(7, 6).Value = 'ДГИ-В-87745/24  14.08.2024'
$ws.Cells.Item(7, 7).Value = 'ДГИ-В-87745/24 14.08.2024'
$ws.Cells.Item(7, 8).Value = 'Кому: Спесивцева С.В. (Департамент городского имущества города Москвы), Мусиенко О.А. (Департамент городского имущества города Москвы)  От кого:  Быкова О.В. (Департамент городского имущества города Москвы)'
$ws.Cells.Item(7, 9).Value = 'УОКУиРП. Ответ. О регистрации права собственности города Москвы на жилой дом по адресу: г. Москва, г. Зеленоград, ул. Заречная, д. 29.'

# Row 8
$ws.Cells.Item(8, 5).Value = 'Вх'
$ws.Cells.Item(8, 6).Value = 'ДГИ-Ф-2863/24  14.08.2024'
$ws.Cells.Item(8, 7).Value = 'ДГП-03-7213/24 14.08.2024'
$ws.Cells.Item(8, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы), Валуй А.А. (Департамент градостроительной политики города Москвы), Караванова Н.П. (Департамент градостроительной политики города Москвы), Курилов А.Ф. (Департамент градостроительной политики города Москвы), Торсунов В.Ю. (Департамент жилищно-коммунального хозяйства города Москвы), Беляев А.А. (Департамент капитального ремонта города Москвы), Загрутдинов Р.Р. (Департамент строительства города Москвы), Жидкин В.Ф. (Департамент развития новых территорий города Москвы), Княжевская Ю.В. (Комитет по архитектуре и градостроительству города Москвы), Щербаков И.А. (Комитет города Москвы по ценовой политике в строительстве и государственной экспертизе проектов), Слободчиков А.О. (Комитет государственного строительного надзора города Москвы)  От кого:  Стулов Д.Ю. (Департамент градостроительной политики города Москвы)'
$ws.Cells.Item(8, 9).Value = 'ФАКСОГРАММА 14 августа 2024 года в 13 часов 00 минут в Департаменте градостроительной политики города Москвы в режиме ВКС состоится совещание по вопросу включения в государственные программы города Москвы «Жилище» и «Градостроительная политика» ключевых показателей, характеризующих эффективность от реализации стратегических задач. Ссылка для участия в совещании: https://vks.dgp.mos.ru/c/70601'

# Row 9
$ws.Cells.Item(9, 5).Value = 'Гр'
$ws.Cells.Item(9, 6).Value = 'ДГИ-ЭГР-46931/24  14.08.2024'
$ws.Cells.Item(9, 7).Value = '56022495 13.08.2024'
$ws.Cells.Item(9, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Обращение граждан (Обращение граждан)'
$ws.Cells.Item(9, 9).Value = 'Обращения граждан Вопрос 1. Сообщение с mos.ru, идентификатор: 56022495 Талызин Александр Викторович, Об уведомлении об улучшении жилищных условий'

# Row 10
$ws.Cells.Item(10, 5).Value = 'Вх'
$ws.Cells.Item(10, 6).Value = 'ДГИ-1-44675/24  13.08.2024'
$ws.Cells.Item(10, 7).Value = '4 02.08.2024'
$ws.Cells.Item(10, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Ермоленко Н.В. (Финансовый управляющий)'
$ws.Cells.Item(10, 9).Value = 'ДГИ-241024/24-(0)-0 запрос по делу А40-95663/2024'

# Row 11
$ws.Cells.Item(11, 5).Value = 'Вх'
$ws.Cells.Item(11, 6).Value = 'ДГИ-1-44670/24  13.08.2024'
$ws.Cells.Item(11, 7).Value = '2 05.08.2024'
$ws.Cells.Item(11, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Мокрушин С.В. (Конкурсный управляющий)'
$ws.Cells.Item(11, 9).Value = 'ДГИ-241090/24-(0)-0 запрос по делу А40-287374/23'

# Row 12
$ws.Cells.Item(12, 5).Value = 'Вх'
$ws.Cells.Item(12, 6).Value = 'ДГИ-1-44664/24  13.08.2024'
$ws.Cells.Item(12, 7).Value = 'Б/Н 02.08.2024'
$ws.Cells.Item(12, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Алехин Н.Н. (Конкурсный управляющий)'
$ws.Cells.Item(12, 9).Value = 'ДГИ-240927/24-(0)-0 запрос по делу А40-169117/23'

# Row 13
$ws.Cells.Item(13, 5).Value = 'Вх'
$ws.Cells.Item(13, 6).Value = 'ДГИ-1-44660/24  13.08.2024'
$ws.Cells.Item(13, 7).Value = '9 02.08.2024'
$ws.Cells.Item(13, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Курзин Д.А. (Финансовый управляющий)'
$ws.Cells.Item(13, 9).Value = 'ДГИ-240954/24-(0)-0 запрос по делу А40-111870/2024'

# Row 14
$ws.Cells.Item(14, 5).Value = 'Вх'
$ws.Cells.Item(14, 6).Value = 'ДГИ-1-44658/24  13.08.2024'
$ws.Cells.Item(14, 7).Value = '1 04.06.2024'
$ws.Cells.Item(14, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Вахрушев В.О. (Временный управляющий)'
$ws.Cells.Item(14, 9).Value = 'ДГИ-240946/24-(0)-0 запрос по делу А40-79798/24'

# Row 15
$ws.Cells.Item(15, 5).Value = 'Вх'
$ws.Cells.Item(15, 6).Value = 'ДГИ-1-44654/24  13.08.2024'
$ws.Cells.Item(15, 7).Value = '9 05.08.2024'
$ws.Cells.Item(15, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Ломакина М.М. (Финансовый управляющий)'
$ws.Cells.Item(15, 9).Value = 'ДГИ-241014/24-(0)-0 запрос по делу А40-73286/2024'

# Row 16
$ws.Cells.Item(16, 5).Value = 'Вх'
$ws.Cells.Item(16, 6).Value = 'ДГИ-1-44650/24  13.08.2024'
$ws.Cells.Item(16, 7).Value = '859-47 07.08.2024'
$ws.Cells.Item(16, 8).Value = 'Кому: Гаман М.Ф. (Департамент городского имущества города Москвы)  От кого:  Стародубцев А.В. (Конкурсный управляющий)'
$ws.Cells.Item(16, 9).Value = 'ДГИ-240940/24-(0)-0 запрос по делу А40-253586/23'
